# Auto-generated edit script applying scheduled market-data refresh
# to the Garuda_Profits sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("H40").Value = 1761.25
$ws.Range("I40").Value = 1626.6666
$ws.Range("J40").Value = 1792.3077
$ws.Range("K40").Value = 1626.6666
$ws.Range("L40").Value = 1792.3077
$ws.Range("M40").Value = -1451.6666
$ws.Range("N40").Value = -2142.3077
# row 76
$ws.Range("H76").Value = 169075.5
$ws.Range("I76").Value = 252213.25
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 252213.25
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -251898.25
$ws.Range("N76").Value = -3430
# row 79
$ws.Range("H79").Value = 169075.5
$ws.Range("I79").Value = 252213.25
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 252213.25
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -251121.25
$ws.Range("N79").Value = -4984
# row 132
$ws.Range("H132").Value = 3041260.8
$ws.Range("I132").Value = 3761142
$ws.Range("K132").Value = 11283426
$ws.Range("M132").Value = -11280896

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3587.95
$ws.Range("I32").Value = 2716.809
$ws.Range("J32").Value = 10636.272
$ws.Range("K32").Value = 2716.809
$ws.Range("L32").Value = 10636.272
$ws.Range("M32").Value = -2429.809
$ws.Range("N32").Value = -11210.272
# row 37
$ws.Range("H37").Value = 9234.666999999999
$ws.Range("I37").Value = 2934
$ws.Range("J37").Value = 12385
$ws.Range("K37").Value = 2934
$ws.Range("L37").Value = 12385
$ws.Range("M37").Value = -2661
$ws.Range("N37").Value = -12931
# row 44
$ws.Range("H44").Value = 21233.334
$ws.Range("J44").Value = 21233.334
$ws.Range("L44").Value = 21233.334
$ws.Range("N44").Value = -22209.334
# row 55
$ws.Range("H55").Value = 22600
$ws.Range("J55").Value = 22600
$ws.Range("L55").Value = 22600
$ws.Range("N55").Value = -23230
# row 63
$ws.Range("H63").Value = 1252024.9
$ws.Range("I63").Value = 1430485.6
$ws.Range("K63").Value = 1430485.6
$ws.Range("M63").Value = -1429799.6
# row 66
$ws.Range("H66").Value = 1252024.9
$ws.Range("I66").Value = 1430485.6
$ws.Range("K66").Value = 7152428
$ws.Range("M66").Value = -7148996

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 404.1905
$ws.Range("I22").Value = 410.44446
$ws.Range("J22").Value = 366.66666
$ws.Range("K22").Value = 410.44446
$ws.Range("L22").Value = 366.66666
$ws.Range("M22").Value = -237.44446
$ws.Range("N22").Value = -712.66666
# row 35
$ws.Range("H35").Value = 29800
$ws.Range("J35").Value = 29800
$ws.Range("L35").Value = 29800
$ws.Range("N35").Value = -30420
# row 82
$ws.Range("H82").Value = 11374.889
$ws.Range("I82").Value = 3427.7144
$ws.Range("J82").Value = 39190
$ws.Range("K82").Value = 3427.7144
$ws.Range("L82").Value = 39190
$ws.Range("M82").Value = -3044.7144
$ws.Range("N82").Value = -39956
# row 85
$ws.Range("H85").Value = 11374.889
$ws.Range("I85").Value = 3427.7144
$ws.Range("J85").Value = 39190
$ws.Range("K85").Value = 3427.7144
$ws.Range("L85").Value = 39190
$ws.Range("M85").Value = -2101.7144
$ws.Range("N85").Value = -41842

$ws = $wb.Worksheets.Item("CRP")
# row 17
$ws.Range("H17").Value = 1900
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1900
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1900
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2248
# row 22
$ws.Range("H22").Value = 373.8125
$ws.Range("I22").Value = 361.72726
$ws.Range("J22").Value = 400.4
$ws.Range("K22").Value = 361.72726
$ws.Range("L22").Value = 400.4
$ws.Range("M22").Value = -11.72726
$ws.Range("N22").Value = -1100.4
# row 25
$ws.Range("H25").Value = 5013
$ws.Range("J25").Value = 5013
$ws.Range("L25").Value = 5013
$ws.Range("N25").Value = -5361
# row 41
$ws.Range("H41").Value = 12043.333
$ws.Range("J41").Value = 12043.333
$ws.Range("L41").Value = 12043.333
$ws.Range("N41").Value = -12899.333
# row 50
$ws.Range("H50").Value = 15410
$ws.Range("J50").Value = 15410
$ws.Range("L50").Value = 15410
$ws.Range("N50").Value = -16660
# row 51
$ws.Range("H51").Value = 21562.8
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 21562.8
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 21562.8
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -23034.8
# row 59
$ws.Range("H59").Value = 29016.25
$ws.Range("I59").Value = 20000
$ws.Range("J59").Value = 30304.285
$ws.Range("K59").Value = 20000
$ws.Range("L59").Value = 30304.285
$ws.Range("M59").Value = -18855
$ws.Range("N59").Value = -32594.285
# row 60
$ws.Range("H60").Value = 18775.6
$ws.Range("I60").Value = 12546.5
$ws.Range("J60").Value = 20332.875
$ws.Range("K60").Value = 12546.5
$ws.Range("L60").Value = 20332.875
$ws.Range("M60").Value = -12035.5
$ws.Range("N60").Value = -21354.875
# row 61
$ws.Range("H61").Value = 21562.8
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 21562.8
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 21562.8
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -22258.8
# row 68
$ws.Range("H68").Value = 38924.75
$ws.Range("J68").Value = 38924.75
$ws.Range("L68").Value = 38924.75
$ws.Range("N68").Value = -40422.75
# row 71
$ws.Range("H71").Value = 38924.75
$ws.Range("J71").Value = 38924.75
$ws.Range("L71").Value = 116774.25
$ws.Range("N71").Value = -124262.25
# row 74
$ws.Range("H74").Value = 25057
$ws.Range("I74").Value = 10285
$ws.Range("J74").Value = 27519
$ws.Range("K74").Value = 10285
$ws.Range("L74").Value = 27519
$ws.Range("M74").Value = -9411
$ws.Range("N74").Value = -29267
# row 77
$ws.Range("H77").Value = 25057
$ws.Range("I77").Value = 10285
$ws.Range("J77").Value = 27519
$ws.Range("K77").Value = 30855
$ws.Range("L77").Value = 82557
$ws.Range("M77").Value = -26487
$ws.Range("N77").Value = -91293
# row 86
$ws.Range("H86").Value = 333335840
$ws.Range("I86").Value = 500001500
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 500001500
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -500000377
$ws.Range("N86").Value = -6746
# row 89
$ws.Range("H89").Value = 333335840
$ws.Range("I89").Value = 500001500
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 2500007500
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -2500001884
$ws.Range("N89").Value = -33732
# row 109
$ws.Range("H109").Value = 9500
$ws.Range("J109").Value = 9500
$ws.Range("L109").Value = 9500
$ws.Range("N109").Value = -11580
# row 132
$ws.Range("H132").Value = 4313561
$ws.Range("I132").Value = 2722.4
$ws.Range("J132").Value = 8932316
$ws.Range("K132").Value = 8167.200000000001
$ws.Range("L132").Value = 26796948
$ws.Range("M132").Value = -5637.200000000001
$ws.Range("N132").Value = -26802008

$ws = $wb.Worksheets.Item("CUL")
# row 6
$ws.Range("H6").Value = 576.2
$ws.Range("I6").Value = 68.25
$ws.Range("J6").Value = 914.8333
$ws.Range("K6").Value = 204.75
$ws.Range("L6").Value = 2744.4999
$ws.Range("M6").Value = -91.75
$ws.Range("N6").Value = -2970.4999
# row 131
$ws.Range("H131").Value = 768.24
$ws.Range("I131").Value = 424
$ws.Range("J131").Value = 786.3579
$ws.Range("K131").Value = 1272
$ws.Range("L131").Value = 2359.0737
$ws.Range("M131").Value = 3768
$ws.Range("N131").Value = -12439.0737

$ws = $wb.Worksheets.Item("GSM")
# row 123
$ws.Range("H123").Value = 29594.428
$ws.Range("J123").Value = 29594.428
$ws.Range("L123").Value = 29594.428
$ws.Range("N123").Value = -34494.428
# row 131
$ws.Range("H131").Value = 22765.428
$ws.Range("J131").Value = 22765.428
$ws.Range("L131").Value = 22765.428
$ws.Range("N131").Value = -32845.428

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 539.625
$ws.Range("I22").Value = 629.1667
$ws.Range("J22").Value = 271
$ws.Range("K22").Value = 629.1667
$ws.Range("L22").Value = 271
$ws.Range("M22").Value = -334.1667
$ws.Range("N22").Value = -861
# row 27
$ws.Range("H27").Value = 539.625
$ws.Range("I27").Value = 629.1667
$ws.Range("J27").Value = 271
$ws.Range("K27").Value = 629.1667
$ws.Range("L27").Value = 271
$ws.Range("M27").Value = -522.1667
$ws.Range("N27").Value = -485
# row 133
$ws.Range("H133").Value = 25708.4
$ws.Range("J133").Value = 25708.4
$ws.Range("L133").Value = 25708.4
$ws.Range("N133").Value = -30768.4

$ws = $wb.Worksheets.Item("WVR")
# row 109
$ws.Range("H109").Value = 12999.333
$ws.Range("J109").Value = 12999.333
$ws.Range("L109").Value = 12999.333
$ws.Range("N109").Value = -15773.333

